$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")
$wsAssets    = $wb.Worksheets.Item("Assets")

# --- Constants sheet: add Retry settings for Get Transaction Item / Set Transaction Status ---
# Write the "Name" cells first, then the "Description" cells, then the numeric
# "Value" cells so new shared-string entries land in the same order as the
# target workbook (Name1, Name2, Desc-row12, Desc1, Desc2).
$wsConstants.Range("A14").Value = "RetryNumberGetTransactionItem"
$wsConstants.Range("A15").Value = "RetryNumberSetTransactionStatus"

$wsConstants.Range("C12").Value = "Error message in case MaxConsecutiveSystemExceptions number is reached."
$wsConstants.Range("C14").Value = "The number of times Get Transaction Item activity is retried in case of an exception. Must be an integer >= 1."
$wsConstants.Range("C15").Value = "The number of times Set transaction status activity is retried in case of an exception. Must be an integer >= 1. "

$wsConstants.Range("B14").Value = 2
$wsConstants.Range("B15").Value = 2

# --- Selection / active-sheet bookkeeping (moved focus from Settings to Constants) ---
[void]$wsSettings.Range("A20").Select()
[void]$wsAssets.Range("B10").Select()
[void]$wsConstants.Activate()
[void]$wsConstants.Range("C22").Select()
